$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Replace the first occurrence of $searchText (found via plain IndexOf on the
# whole document text) with $replaceText. Returns the absolute start offset
# of the freshly-inserted replacement text (so callers can re-split it).
function Replace-Segment {
    param(
        [string]$searchText,
        [string]$replaceText
    )
    $full = $d.Content.Text
    $idx = $full.IndexOf($searchText)
    if ($idx -lt 0) {
        throw "Replace-Segment: text not found: $searchText"
    }
    $len = $searchText.Length
    $r = $d.Range($idx, $idx + $len)
    $r.Text = $replaceText
    return $idx
}

# Force a run-split boundary at absolute offset $at (within the still-open
# segment that ends at absolute offset $segEnd) without altering the visible
# formatting: toggle Bold on then immediately back to "undefined" so the
# rPr serializes identically to its neighbours but the run list splits.
function Split-At {
    param(
        [int]$at,
        [int]$segEnd
    )
    $r = $d.Range($at, $segEnd)
    $r.Bold = 1
    $r.Bold = 9999999
}

# Replace a whole segment (text between hard breaks / paragraph boundaries)
# with a sequence of pieces, re-splitting the resulting single run back into
# one run per piece so the final run layout has one run per $pieces element.
function Replace-SegmentPieces {
    param(
        [string]$searchText,
        [string[]]$pieces
    )
    $replaceText = [string]::Join("", $pieces)
    $base = Replace-Segment $searchText $replaceText
    $segEnd = $base + $replaceText.Length

    $offset = 0
    $splitPoints = New-Object System.Collections.Generic.List[int]
    for ($i = 0; $i -lt $pieces.Length - 1; $i++) {
        $offset += $pieces[$i].Length
        $splitPoints.Add($offset)
    }
    foreach ($pos in $splitPoints) {
        Split-At ($base + $pos) $segEnd
    }
}

# ---------------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------------
Replace-Segment "Unraveling the Enigmatic World of Quantum Computing" `
    "The Intricate Connection: How Politics, History, and Government Interweave to Shape Our World" | Out-Null

# ---------------------------------------------------------------------------
# Author name ("Dr" + "." + " Emily White" -> single run "Sophia Harris")
# ---------------------------------------------------------------------------
Replace-Segment "Dr. Emily White" "Sophia Harris" | Out-Null

# ---------------------------------------------------------------------------
# Email ("emily" + "." + "white@quantumresearch" -> single run;
# the trailing "." + "com" runs are left untouched)
# ---------------------------------------------------------------------------
Replace-Segment "emily.white@quantumresearch" "sophiaharris@validemail" | Out-Null

# ---------------------------------------------------------------------------
# Body paragraph, segment 1 (before the first line break)
# ---------------------------------------------------------------------------
Replace-SegmentPieces `
    "In the vast and ever-evolving realm of computer science, a new frontier has emerged, beckoning us to explore the enigmatic world of quantum computing. This revolutionary field promises to transform our understanding of computation and propel us into an era of unprecedented technological advancement. Quantum computers, harnessing the power of quantum mechanics, operate on principles vastly different from their classical counterparts, offering the potential to solve complex problems that have remained intractable for conventional computers." `
    @(
        "In the tapestry of human society, politics, history, and government intertwine to form an intricate web that profoundly shapes our world",
        ".",
        " Politics, as the art of governance, involves the allocation of power and resources within a society",
        ".",
        " History, as the study of past events, provides insights into how societies have evolved and how political decisions have impacted them",
        ".",
        " Government, as the system through which power is exercised, establishes rules and regulations that govern the behavior of individuals and organizations",
        ".",
        " These three spheres are inseparable, influencing and influencing each other in myriad ways",
        "."
    )

# ---------------------------------------------------------------------------
# Body paragraph, segment 2 (between the two line-break pairs)
# ---------------------------------------------------------------------------
Replace-SegmentPieces `
    "Delving into the intricacies of quantum computing unveils a fascinating tapestry of concepts that challenge our traditional notions of computation. This paradigm shift involves the manipulation of quantum bits, or qubits, which exist in a superposition of states, enabling them to encode information in a manner that classical bits cannot. Furthermore, quantum mechanics introduces the phenomenon of entanglement, where the state of one qubit becomes instantaneously correlated with the state of another, regardless of the distance between them. These remarkable properties empower quantum computers to process vast amounts of data concurrently, tackling computational challenges that were previously deemed insurmountable." `
    @(
        "Comprehending the relationship between politics, history, and government is crucial for understanding how societies function",
        ".",
        " Politics is often a reflection of the power dynamics within a society, with different groups competing for influence and control",
        ".",
        " Historical events can shape political ideologies and institutions, while government policies can have profound historical consequences",
        ".",
        " Effective governance requires an examination of both political and historical factors, ensuring policies are informed by past lessons and adapted to current realities",
        "."
    )

# ---------------------------------------------------------------------------
# Body paragraph, segment 3 (after the second line-break pair)
# ---------------------------------------------------------------------------
Replace-SegmentPieces `
    "As quantum computing matures, its potential applications span a broad spectrum of fields, poised to revolutionize industries and reshape our world. From unraveling the mysteries of protein folding to designing novel materials with exceptional properties, quantum computers hold the promise of unlocking breakthroughs in fields ranging from medicine and finance to cryptography and beyond. This nascent technology has the potential to reshape the very fabric of our digital infrastructure, ushering in an era of enhanced security, accelerated simulations, and unprecedented computational power." `
    @(
        "The interconnectedness of politics, history, and government is evident throughout history",
        ".",
        " The American Revolution, for instance, was both a political struggle for independence and a pivotal historical event that shaped the course of American history",
        ".",
        " The subsequent establishment of the United States government was a direct outcome of this political and historical upheaval",
        ".",
        " Similarly, the rise of communism in the 20th century had profound political, historical, and governmental implications worldwide",
        "."
    )

# ---------------------------------------------------------------------------
# Summary heading stays "Summary" (untouched)
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Summary body paragraph (one run per sentence already; same run count)
# ---------------------------------------------------------------------------
Replace-Segment "Quantum computing, a paradigm-shifting field at the forefront of computer science, harnesses the principles of quantum mechanics to unlock unprecedented computational power" `
    "Politics, history, and government are inextricably linked, shaping the fabric of human society" | Out-Null

Replace-Segment " Quantum " `
    " Politics reflects power dynamics and decision-making, history provides context for " | Out-Null

Replace-Segment "bits, or qubits, possess unique properties like superposition and entanglement, enabling them to process vast amounts of data concurrently and tackle problems that are intractable for classical computers" `
    "understanding contemporary issues, and government establishes rules and regulations" | Out-Null

Replace-Segment " The potential applications of quantum computing are vast, with implications for fields such as medicine, finance, cryptography, and materials science" `
    " Their interplay influences the direction of societies, from shaping political ideologies to determining the distribution of resources" | Out-Null

Replace-Segment " As research continues to advance, quantum computing holds the promise of revolutionizing industries and propelling us into a new era of technological innovation" `
    " Comprehending this interconnectedness is essential for informed citizenship and effective governance, allowing us to navigate the complexities of our world with greater understanding and agency" | Out-Null

# ---------------------------------------------------------------------------
# Add a new empty paragraph at the very end of the document body
# ---------------------------------------------------------------------------
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()
